$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# ---------------------------------------------------------------------------
# Row 2 / Row 3 / Row 6 keep their existing semantic formatting; nothing to
# touch there directly. The only real formatting change in that block is
# that row 3 (A3:I3) gets "wrap text" turned on and a taller row height.
# ---------------------------------------------------------------------------
$ws.Range("A3:I3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 61.5

# ---------------------------------------------------------------------------
# Row 8 content updates (2021 4th-quarter report -> 2022 1st-quarter report)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value2 = 2022
$ws.Range("B8").Value2 = 44562
$ws.Range("C8").Value2 = 44651

# J8 text gains a leading space: "Secretaria Academica (UPP)" -> " Secretaria Academica (UPP)"
$ws.Range("J8").Value2 = " Secretaria Academica (UPP)"

$ws.Range("K8").Value2 = 44659
$ws.Range("L8").Value2 = 44659

# M8 keeps its "Conforme..." note text - untouched.

$ws.Rows.Item(8).RowHeight = 105

# ---------------------------------------------------------------------------
# Column M width shrinks slightly
# ---------------------------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = 83.67

# ---------------------------------------------------------------------------
# Data validation on column E only applies to E8 now (was E8:E201)
# ---------------------------------------------------------------------------
$ws.Range("E9:E201").Validation.Delete()

# ---------------------------------------------------------------------------
# Sheet view: scroll back to top-left and move the selection to B14
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("B14").Select()
